$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1154
$ws.Range("I12").Value = 644.75
$ws.Range("J12").Value = 1833
$ws.Range("K12").Value = 644.75
$ws.Range("L12").Value = 1833
$ws.Range("M12").Value = -474.75
$ws.Range("N12").Value = -2173
$ws.Range("H17").Value = 812710.6
$ws.Range("J17").Value = 858657.8
$ws.Range("L17").Value = 2575973.4
$ws.Range("N17").Value = -2576309.4
$ws.Range("H19").Value = 1404.6316
$ws.Range("I19").Value = 818.8
$ws.Range("J19").Value = 1613.8572
$ws.Range("K19").Value = 818.8
$ws.Range("L19").Value = 1613.8572
$ws.Range("M19").Value = -643.8
$ws.Range("N19").Value = -1963.8572
$ws.Range("H32").Value = 3619.8
$ws.Range("J32").Value = 3619.8
$ws.Range("L32").Value = 3619.8
$ws.Range("N32").Value = -4271.8
$ws.Range("H70").Value = 4173
$ws.Range("I70").Value = 5226.3335
$ws.Range("J70").Value = 2593
$ws.Range("K70").Value = 15679.0005
$ws.Range("L70").Value = 7779
$ws.Range("M70").Value = -15409.0005
$ws.Range("N70").Value = -8319
$ws.Range("H73").Value = 4173
$ws.Range("I73").Value = 5226.3335
$ws.Range("J73").Value = 2593
$ws.Range("K73").Value = 15679.0005
$ws.Range("L73").Value = 7779
$ws.Range("M73").Value = -14743.0005
$ws.Range("N73").Value = -9651
$ws.Range("H113").Value = 9095
$ws.Range("I113").Value = 9738.571
$ws.Range("K113").Value = 9738.571
$ws.Range("M113").Value = -6484.571
$ws.Range("H116").Value = 859981.4399999999
$ws.Range("I116").Value = 2782327.2
$ws.Range("J116").Value = 5605.5557
$ws.Range("K116").Value = 2782327.2
$ws.Range("L116").Value = 5605.5557
$ws.Range("M116").Value = -2778885.2
$ws.Range("N116").Value = -12489.5557
$ws.Range("H137").Value = 18632.385
$ws.Range("I137").Value = 21111.455
$ws.Range("K137").Value = 63334.36500000001
$ws.Range("M137").Value = -60784.36500000001
$ws.Range("H138").Value = 3991.1892
$ws.Range("J138").Value = 4600.517
$ws.Range("L138").Value = 13801.551
$ws.Range("N138").Value = -24081.551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14524.8
$ws.Range("I2").Value = 24601.6
$ws.Range("K2").Value = 24601.6
$ws.Range("M2").Value = -24488.6
$ws.Range("H5").Value = 388
$ws.Range("I5").Value = 417
$ws.Range("J5").Value = 337.25
$ws.Range("K5").Value = 417
$ws.Range("L5").Value = 337.25
$ws.Range("M5").Value = -305
$ws.Range("N5").Value = -561.25
$ws.Range("H43").Value = 47499.5
$ws.Range("J43").Value = 47499.5
$ws.Range("L43").Value = 47499.5
$ws.Range("N43").Value = -48125.5
$ws.Range("H61").Value = 3886.1
$ws.Range("I61").Value = 3112.8235
$ws.Range("J61").Value = 5529.3125
$ws.Range("K61").Value = 3112.8235
$ws.Range("L61").Value = 5529.3125
$ws.Range("M61").Value = -2900.8235
$ws.Range("N61").Value = -5953.3125
$ws.Range("H97").Value = 50857.285
$ws.Range("I97").Value = 28747.5
$ws.Range("J97").Value = 80337
$ws.Range("K97").Value = 28747.5
$ws.Range("L97").Value = 80337
$ws.Range("M97").Value = -28251.5
$ws.Range("N97").Value = -81329
$ws.Range("H116").Value = 14524.8
$ws.Range("I116").Value = 24601.6
$ws.Range("K116").Value = 24601.6
$ws.Range("M116").Value = -22307.6
$ws.Range("H132").Value = 3575.3794
$ws.Range("I132").Value = 1814.625
$ws.Range("J132").Value = 5742.4614
$ws.Range("K132").Value = 5443.875
$ws.Range("L132").Value = 17227.3842
$ws.Range("M132").Value = -2913.875
$ws.Range("N132").Value = -22287.3842
$ws.Range("H133").Value = 64999.668
$ws.Range("J133").Value = 64999.668
$ws.Range("L133").Value = 64999.668
$ws.Range("N133").Value = -70059.66800000001
$ws.Range("H136").Value = 3886.1
$ws.Range("I136").Value = 3112.8235
$ws.Range("J136").Value = 5529.3125
$ws.Range("K136").Value = 9338.470499999999
$ws.Range("L136").Value = 16587.9375
$ws.Range("M136").Value = -6788.470499999999
$ws.Range("N136").Value = -21687.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14524.8
$ws.Range("I3").Value = 24601.6
$ws.Range("K3").Value = 24601.6
$ws.Range("M3").Value = -24487.6
$ws.Range("H4").Value = 388
$ws.Range("I4").Value = 417
$ws.Range("J4").Value = 337.25
$ws.Range("K4").Value = 417
$ws.Range("L4").Value = 337.25
$ws.Range("M4").Value = -302
$ws.Range("N4").Value = -567.25
$ws.Range("H20").Value = 4057.3
$ws.Range("I20").Value = 1217.6
$ws.Range("J20").Value = 6897
$ws.Range("K20").Value = 1217.6
$ws.Range("L20").Value = 6897
$ws.Range("M20").Value = -970.5999999999999
$ws.Range("N20").Value = -7391
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -711

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 691.8182
$ws.Range("I22").Value = 605.5
$ws.Range("K22").Value = 605.5
$ws.Range("M22").Value = -255.5
$ws.Range("H31").Value = 8328.833000000001
$ws.Range("I31").Value = 11987.5
$ws.Range("J31").Value = 6499.5
$ws.Range("K31").Value = 11987.5
$ws.Range("L31").Value = 6499.5
$ws.Range("M31").Value = -11692.5
$ws.Range("N31").Value = -7089.5
$ws.Range("H34").Value = 8328.833000000001
$ws.Range("I34").Value = 11987.5
$ws.Range("J34").Value = 6499.5
$ws.Range("K34").Value = 11987.5
$ws.Range("L34").Value = 6499.5
$ws.Range("M34").Value = -11785.5
$ws.Range("N34").Value = -6903.5
$ws.Range("H58").Value = 2229.5881
$ws.Range("I58").Value = 1619.6154
$ws.Range("K58").Value = 1619.6154
$ws.Range("M58").Value = -1416.6154
$ws.Range("H59").Value = 74999.5
$ws.Range("J59").Value = 74999.5
$ws.Range("L59").Value = 74999.5
$ws.Range("N59").Value = -77289.5
$ws.Range("H136").Value = 2229.5881
$ws.Range("I136").Value = 1619.6154
$ws.Range("K136").Value = 4858.8462
$ws.Range("M136").Value = -2308.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 298.35483
$ws.Range("I26").Value = 292.56
$ws.Range("K26").Value = 877.6800000000001
$ws.Range("M26").Value = -589.6800000000001
$ws.Range("H51").Value = 2050.8948
$ws.Range("J51").Value = 2666.3333
$ws.Range("L51").Value = 7998.999899999999
$ws.Range("N51").Value = -8918.999899999999
$ws.Range("H63").Value = 2622
$ws.Range("J63").Value = 2996
$ws.Range("L63").Value = 8988
$ws.Range("N63").Value = -10486
$ws.Range("H66").Value = 2622
$ws.Range("J66").Value = 2996
$ws.Range("L66").Value = 26964
$ws.Range("N66").Value = -34452

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20312
$ws.Range("H70").Value = 7325.846
$ws.Range("I70").Value = 7208
$ws.Range("J70").Value = 7378.222
$ws.Range("K70").Value = 7208
$ws.Range("L70").Value = 7378.222
$ws.Range("M70").Value = -6938
$ws.Range("N70").Value = -7918.222
$ws.Range("H73").Value = 7325.846
$ws.Range("I73").Value = 7208
$ws.Range("J73").Value = 7378.222
$ws.Range("K73").Value = 7208
$ws.Range("L73").Value = 7378.222
$ws.Range("M73").Value = -6272
$ws.Range("N73").Value = -9250.222
$ws.Range("H80").Value = 3999.5
$ws.Range("J80").Value = 3999.5
$ws.Range("L80").Value = 3999.5
$ws.Range("N80").Value = -5995.5
$ws.Range("H83").Value = 3999.5
$ws.Range("J83").Value = 3999.5
$ws.Range("L83").Value = 19997.5
$ws.Range("N83").Value = -29981.5
$ws.Range("H102").Value = 7975.143
$ws.Range("I102").Value = 9654.875
$ws.Range("K102").Value = 9654.875
$ws.Range("M102").Value = -8032.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3989.35
$ws.Range("I16").Value = 3321.4443
$ws.Range("K16").Value = 3321.4443
$ws.Range("M16").Value = -3151.4443
$ws.Range("H55").Value = 1280.5
$ws.Range("I55").Value = 262.33334
$ws.Range("J55").Value = 1789.5834
$ws.Range("K55").Value = 262.33334
$ws.Range("L55").Value = 1789.5834
$ws.Range("M55").Value = -89.33334000000002
$ws.Range("N55").Value = -2135.5834
$ws.Range("H64").Value = 29716.666
$ws.Range("J64").Value = 29716.666
$ws.Range("L64").Value = 29716.666
$ws.Range("N64").Value = -30166.666
$ws.Range("H67").Value = 29716.666
$ws.Range("J67").Value = 29716.666
$ws.Range("L67").Value = 29716.666
$ws.Range("N67").Value = -31276.666
$ws.Range("H136").Value = 4815.436
$ws.Range("I136").Value = 1755.8096
$ws.Range("K136").Value = 5267.4288
$ws.Range("M136").Value = -2717.4288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2608.1304
$ws.Range("I136").Value = 1144.0714
$ws.Range("J136").Value = 4885.5557
$ws.Range("K136").Value = 3432.2142
$ws.Range("L136").Value = 14656.6671
$ws.Range("M136").Value = -882.2142000000003
$ws.Range("N136").Value = -19756.6671
